# "perbaikan brian 26 januari 2025"
# Re-work the barang/toko/gudang import template:
#  - drop the obsolete leading "status" column
#  - relabel the remaining headers with human-friendly titles
#  - add a sample data row so the import format is self-documenting
#  - resize the columns to fit the new header/content widths

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old layout was: status | nama_barang | nama_gudang | jumlah_barang | nama_toko_luar
# Deleting column A shifts everything left, leaving the 4 columns we still need.
$ws.Range("A1").EntireColumn.Delete()

# Re-label the (now-shifted) headers with proper capitalisation/spacing.
$ws.Range("A1").Value = "Nama Barang"
$ws.Range("B1").Value = "Nama Gudang"
$ws.Range("C1").Value = "Jumlah Barang"
$ws.Range("D1").Value = "Nama Toko Luar"

# Add a sample row illustrating how the import file should be filled in.
$ws.Range("A2").Value = "Travo 5A O"
$ws.Range("B2").Value = "Sumber Agung Gudang"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = "Shopee"

# Widen the columns (character units) to comfortably fit the new content.
$ws.Columns.Item(1).ColumnWidth = 17.451822916666668
$ws.Columns.Item(2).ColumnWidth = 24.877604166666668
$ws.Columns.Item(3).ColumnWidth = 17.877604166666668
$ws.Columns.Item(4).ColumnWidth = 18.022135416666668

# Leave the selection where it was left in the saved file.
$ws.Range("D6").Select() | Out-Null
